# Trade #7 closed at 2026-02-16 21:51:38 - leadlag DOWN +0.000%
# Append the new trade row (row 8) to both the "All Trades" and "leadlag"
# worksheets.

$wb = $excel.ActiveWorkbook

$targetSheets = @("All Trades", "leadlag")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 8

    # --- Plain values (these never get misinterpreted by type inference) ---
    $ws.Cells.Item($row, 1).Value  = 7                      # Trade #
    $ws.Cells.Item($row, 3).Value  = "21:51:38"              # Time
    $ws.Cells.Item($row, 4).Value  = "leadlag"                # Strategy
    $ws.Cells.Item($row, 5).Value  = "DOWN"                    # Side
    $ws.Cells.Item($row, 6).Value  = 68235.28999999999         # Entry Price
    $ws.Cells.Item($row, 8).Value  = "OPEN"                     # Status
    $ws.Cells.Item($row, 9).Value  = 0                           # P&L %
    $ws.Cells.Item($row, 10).Value = 0                            # P&L $
    $ws.Cells.Item($row, 11).Value = 100                           # Capital After
    $ws.Cells.Item($row, 12).Value = 0.75                           # Confidence
    $ws.Cells.Item($row, 13).Value = "Binance leading with -0.098% move" # Entry Reason
    $ws.Cells.Item($row, 15).Value = 0                               # Duration (min)

    # --- Date-like text ("2026-02-16") ---
    # Writing this string straight into Value/Formula makes Excel
    # auto-recognize it as a date and silently reformat the cell. To keep
    # it as plain text (matching the rest of the column), build it as a
    # text formula in a scratch cell, copy it, and paste-special only the
    # value back into the target cell - this bypasses Excel's "smart" input
    # parsing entirely.
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = "=""2026-02-16"""
    $scratch.Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Value = $null
}
